$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at the top; all existing data rows shift down by one.
$ws.Rows("1:1").Insert()

# 2. New header row for the additional measurement columns.
$ws.Range("B1").Value = "Singleton"
$ws.Range("C1").Value = "Transient"
$ws.Range("D1").Value = "Combined"
$ws.Range("E1").Value = "Interception"

# 3. New "Interception" column data (column E) for the rows that have it.
$ws.Range("E2").Value = 75
$ws.Range("E3").Value = 43222

$ws.Range("E18").NumberFormat = "0"
$ws.Range("E18").Value = 560

$ws.Range("E20").NumberFormat = "0"
$ws.Range("E20").Value = 897

$ws.Range("E21").Value = 11603
$ws.Range("E23").Value = 140276
$ws.Range("E24").Value = 26963

# 4. Update the chart series source references so they keep pointing at the
#    same (now shifted-by-one-row) data.
$chart1 = $ws.ChartObjects(1).Chart
$s1 = $chart1.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(,Tabelle1!`$A`$3:`$A`$24,Tabelle1!`$B`$3:`$B`$24,1)"

$chart2 = $ws.ChartObjects(2).Chart
$s2 = $chart2.SeriesCollection().Item(1)
$s2.Formula = "=SERIES(,(Tabelle1!`$A`$3:`$A`$18,Tabelle1!`$A`$20,Tabelle1!`$A`$21:`$A`$24),(Tabelle1!`$C`$3:`$C`$18,Tabelle1!`$C`$20,Tabelle1!`$C`$21:`$C`$24),1)"

$chart3 = $ws.ChartObjects(3).Chart
$s3 = $chart3.SeriesCollection().Item(1)
$s3.Formula = "=SERIES(,(Tabelle1!`$A`$3:`$A`$18,Tabelle1!`$A`$20,Tabelle1!`$A`$21:`$A`$24),(Tabelle1!`$D`$3:`$D`$18,Tabelle1!`$D`$20,Tabelle1!`$D`$21:`$D`$24),1)"

# 5. Move the three chart objects down by one row to stay anchored where
#    they were relative to the (now shifted) data / empty rows above them.
$co1 = $ws.ChartObjects(1)
$co1.Top = $co1.Top + 15

$co2 = $ws.ChartObjects(2)
$co2.Top = $co2.Top + 15

$co3 = $ws.ChartObjects(3)
$co3.Top = $co3.Top + 15
$co3.Height = $co3.Height + 0.0001

# 6. Restore the active selection.
$ws.Range("E4").Select()
